$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'27.105.39"
$ws.Range("E2").Value = "  -1.64%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.780.01"
$ws.Range("E3").Value = "  -2.17%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.36%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'336.40"

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.22%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.3808"
$ws.Range("E7").Value = "  -0.50%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3406"
$ws.Range("E8").Value = "  -3.29%  "

# Row 9 - OKB
$ws.Range("D9").Value = "'48.00"
$ws.Range("E9").Value = "  -3.58%  "

# Row 10 - Polygon
$ws.Range("D10").Value = "'1.185"
$ws.Range("E10").Value = "  -4.23%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.07424"
$ws.Range("E11").Value = "  -4.85%  "

# Row 12 - BinanceUSD
$ws.Range("D12").Value = "'1.004"
$ws.Range("E12").Value = "  +0.25%  "

# Row 13 - Solana
$ws.Range("D13").Value = "'21.58"
$ws.Range("E13").Value = "  -3.24%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'6.410"
$ws.Range("E14").Value = "  -3.25%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "'1.779.94"
$ws.Range("E15").Value = "  -1.96%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "'7.041"
$ws.Range("E16").Value = "  -2.64%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "'0.00001082"
$ws.Range("E17").Value = "  -4.05%  "

# Row 18 - TRON
$ws.Range("D18").Value = "'0.06643"
$ws.Range("E18").Value = "  -1.40%  "

# Row 19 - Litecoin
$ws.Range("D19").Value = "'83.19"

# Row 20 - Dai
$ws.Range("D20").Value = "'1.003"
$ws.Range("E20").Value = "  +0.24%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'6.540"
$ws.Range("E21").Value = "  +0.04%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  -2.71%  "

# Row 23 - WrappedBTC
$ws.Range("D23").Value = "'27.115.29"
$ws.Range("E23").Value = "  -1.56%  "

# Row 24 - Cosmos
$ws.Range("D24").Value = "'12.19"
$ws.Range("E24").Value = "  -7.90%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "'2.373"
$ws.Range("E25").Value = "  -3.38%  "

# Row 26 - LidoDAOToken
$ws.Range("D26").Value = "'2.500"
$ws.Range("E26").Value = "  -6.88%  "

# Row 27 - ImmutableX
$ws.Range("D27").Value = "'1.456"
$ws.Range("E27").Value = "  -3.17%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'21.03"
$ws.Range("E28").Value = "  -5.24%  "

# Row 29 - Monero
$ws.Range("E29").Value = "  +0.78%  "

# Row 30 - WrappedliquidstakedEther2.0
$ws.Range("D30").Value = "'1.979.26"
$ws.Range("E30").Value = "  -1.95%  "

# Row 31 - BitcoinCash
$ws.Range("D31").Value = "'133.92"
$ws.Range("E31").Value = "  -1.74%  "

# Row 32 - HuobiToken
$ws.Range("D32").Value = "'3.980"
$ws.Range("E32").Value = "  -2.47%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "'6.003"
$ws.Range("E33").Value = "  -6.04%  "

# Row 34 - Stellar
$ws.Range("D34").Value = "'0.08655"
$ws.Range("E34").Value = "  -1.39%  "

# Row 35 - Aptos
$ws.Range("D35").Value = "'13.03"
$ws.Range("E35").Value = "  -7.21%  "

# Row 36 - WEMIXTOKEN
$ws.Range("D36").Value = "'1.620"
$ws.Range("E36").Value = "  -5.13%  "

# Row 37 - InternetComputer(DFINITY)
$ws.Range("D37").Value = "'5.377"
$ws.Range("E37").Value = "  -4.62%  "

# Row 38 - TheSandbox
$ws.Range("E38").Value = "  -4.08%  "

# Row 39 - Hedera
$ws.Range("D39").Value = "'0.06266"
$ws.Range("E39").Value = "  -4.29%  "

# Row 40 / 41 - swap Algorand <-> VeChain with updated values
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.02315"
$ws.Range("E40").Value = "  -4.62%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.2168"
$ws.Range("E41").Value = "  -5.06%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "'8.549"
$ws.Range("E42").Value = "  -5.30%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").Value = "'1.227"
$ws.Range("E43").Value = "  -4.86%  "

# Row 44 / 45 - swap Frax <-> EnergySwap with updated values
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'14.24"
$ws.Range("E44").Value = "  -4.05%  "

$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.003"
$ws.Range("E45").Value = "  +0.30%  "

# Row 46 - Decentraland
$ws.Range("D46").Value = "'0.6395"
$ws.Range("E46").Value = "  -3.57%  "

# Row 47 - PancakeSwap
$ws.Range("D47").Value = "'3.856"
$ws.Range("E47").Value = "  -2.87%  "

# Row 48 - NEARProtocol
$ws.Range("D48").Value = "'2.116"
$ws.Range("E48").Value = "  -3.14%  "

# Row 49 - Quant
$ws.Range("D49").Value = "'131.07"
$ws.Range("E49").Value = "  -1.29%  "

# Row 50 - Cronos
$ws.Range("D50").Value = "'0.07093"

# Row 51 - Aave
$ws.Range("D51").Value = "'78.51"
$ws.Range("E51").Value = "  -2.74%  "
